# Registro de proyectos - agrega columna "Telefono del Coordinador" (K)
# replicando el valor que ya existe en la columna I (comparten el mismo
# texto en la tabla de cadenas compartidas), ajusta el ancho de las
# columnas J y K, y actualiza la celda/ventana seleccionada.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nueva columna K1: mismo texto que I1 ("Teléfono del Coordinador")
$ws.Range("K1").Value = "Teléfono del Coordinador"

# Nuevos anchos de columna (J y K): valores elegidos para que, tras la
# conversion interna caracteres->pixeles de Excel, el ancho guardado en
# el XML quede lo mas cerca posible de 22.28515625 y 19.5703125
$ws.Columns.Item(10).ColumnWidth = 21.451822916666668
$ws.Columns.Item(11).ColumnWidth = 18.736979166666668

# Desplaza la vista (antes topLeftCell = D1, ahora E1) y actualiza la
# celda activa / selección (antes D7, ahora J3)
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J3").Select() | Out-Null
